# aggiornamento excel retail e test XRA_factoring
#
# - Insert a new "EXPOSURE" variable row into the "r AnalysisUnit_Variable"
#   sheet, right after the "SEGMENT" row (new row 6), pushing the IND_*
#   rows down by one.
# - Switch the active tab from "r AnalysisUnit_Variable" back to
#   "Analysis_Unit" (and move the live selection on the variable sheet to
#   F13 before leaving it).
# - Disable concurrent/multi-threaded calculation for the workbook.

$wb = $excel.ActiveWorkbook

$wsAnalysisUnit = $wb.Worksheets.Item("Analysis_Unit")
$wsVariable     = $wb.Worksheets.Item("r AnalysisUnit_Variable")

# --- Insert the new COUNTERPARTY_RETAIL_EXPOSURE / RETAIL_EXPOSURE row ---
$wsVariable.Rows.Item(6).Insert()

$wsVariable.Range("A6").Value = "CREATE/MODIFY"
$wsVariable.Range("B6").Value = "COUNTERPARTY_RETAIL_EXPOSURE"
$wsVariable.Range("C6").Value = "COUNTERPARTY_RETAIL_EXPOSURE"
$wsVariable.Range("E6").Value = "COUNTERPARTY_RETAIL"
$wsVariable.Range("F6").Value = "RETAIL_EXPOSURE"

# --- Move the selection on the variable sheet, then switch back to the
#     Analysis_Unit tab so it becomes the active/selected sheet again ---
$wsVariable.Range("F13").Select()
$wsAnalysisUnit.Select()

# --- Turn off concurrent calculation (calcPr concurrentCalc="0") ---
$excel.MultiThreadedCalculation.Enabled = $false
